$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.030.69"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.270.23"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.66%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.28"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.37"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.31"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.23"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D14").Value = "2.618.91"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").Value = "2.268.29"
$ws.Range("E15").Value = "  -1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.77"
$ws.Range("E16").Value = "  +0.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.801"
$ws.Range("E17").Value = "  -4.86%  "
$ws.Range("D18").Value = "44.885.31"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").Value = "  +6.31%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  -3.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("E21").Value = "  -3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.51"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.56"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.76"
$ws.Range("E27").Value = "  +11.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.56"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.56"
$ws.Range("E30").Value = "  -2.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.27"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -7.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0789"
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.57"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.96"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -5.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.85"
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0307"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.79"
$ws.Range("E42").Value = "  -7.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +12.42%  "
$ws.Range("D45").Value = "1.768.27"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.196"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "70.34"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.10"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.63"
$ws.Range("E49").Value = "  -5.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.85"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.25"
$ws.Range("E51").Value = "  -3.31%  "
